$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'crop athletic leggings'
$ws.Cells.Item(2, 1).Value = 'kid basketball knee pads'
$ws.Cells.Item(3, 1).Value = 'men''s spandex leggings'
$ws.Cells.Item(4, 1).Value = 'hex pants'
$ws.Cells.Item(5, 1).Value = 'medical compression pants'
$ws.Cells.Item(6, 1).Value = 'fitness tights for men'
$ws.Cells.Item(7, 1).Value = 'the rock mens basketball'
$ws.Cells.Item(8, 1).Value = 'black men tights'
$ws.Cells.Item(9, 1).Value = 'd man basketball'
$ws.Cells.Item(10, 1).Value = 'asics compression pants'
$ws.Cells.Item(11, 1).Value = 'knee armor knee pads'
$ws.Cells.Item(12, 1).Value = 'adidas tights for men'
$ws.Cells.Item(13, 1).Value = 'slide on knee pads'
$ws.Cells.Item(14, 1).Value = 'men athletic tights'
$ws.Cells.Item(15, 1).Value = 'blue mens compression pants'
$ws.Cells.Item(16, 1).Value = 'sport tights for men'
$ws.Cells.Item(17, 1).Value = 'compression basketball tights'
$ws.Cells.Item(18, 1).Value = 'kids compression pants'
$ws.Cells.Item(19, 1).Value = 'usa tights men'
$ws.Cells.Item(20, 1).Value = 'razor knee pads'
$ws.Cells.Item(21, 1).Value = 'knee pads addidas'
$ws.Cells.Item(22, 1).Value = 'knee pads hunting'
$ws.Cells.Item(23, 1).Value = 'knee pads leggings'
$ws.Cells.Item(24, 1).Value = 'knee pads elbow pads youth'
$ws.Cells.Item(25, 1).Value = 'knee pads slim'
$ws.Cells.Item(26, 1).Value = 'basketball pants adidas'
$ws.Cells.Item(27, 1).Value = 'mens tights navy'
$ws.Cells.Item(28, 1).Value = 'mens adidas basketball pants'
$ws.Cells.Item(29, 1).Value = 'compression pants men 3xl'
$ws.Cells.Item(30, 1).Value = 'compression pants 2xu'
$ws.Cells.Item(31, 1).Value = 'pro x knee pad'
$ws.Cells.Item(32, 1).Value = 'puma compression pants men'
$ws.Cells.Item(33, 1).Value = 'men''s basketball pants'
$ws.Cells.Item(34, 1).Value = 'gray baseball pants youth girls'
$ws.Cells.Item(35, 1).Value = 'mens basketball jacket'
$ws.Cells.Item(36, 1).Value = 'photography knee pads'
$ws.Cells.Item(37, 1).Value = 'men''s tights leggings'
$ws.Cells.Item(38, 1).Value = 'baseball pants men grey'
$ws.Cells.Item(39, 1).Value = 'venom compression pants'
$ws.Cells.Item(40, 1).Value = 'padded knee tights'
$ws.Cells.Item(41, 1).Value = 'ua compression pants'
$ws.Cells.Item(42, 1).Value = 'men workout tights'
$ws.Cells.Item(43, 1).Value = 'flag compression pants'
$ws.Cells.Item(44, 1).Value = 'cool knee pads'
$ws.Cells.Item(45, 1).Value = 'navy compression leggings'
$ws.Cells.Item(46, 1).Value = 'force knee pads'
$ws.Cells.Item(47, 1).Value = 'mens wrestling pants'
$ws.Cells.Item(48, 1).Value = 'mens capri compression pants'
$ws.Cells.Item(49, 1).Value = 'woman compression pants'
$ws.Cells.Item(50, 1).Value = 'purple knee pads'
$ws.Cells.Item(51, 1).Value = 'reebok knee pads'
$ws.Cells.Item(52, 1).Value = 'venum compression pants men'
$ws.Cells.Item(53, 1).Value = 'purple athletic leggings'
$ws.Cells.Item(54, 1).Value = 'thermal compression pants'
$ws.Cells.Item(55, 1).Value = 'addidas knee pads'
$ws.Cells.Item(56, 1).Value = 'jordan mens tights'
$ws.Cells.Item(57, 1).Value = 'purple compression pants men'
$ws.Cells.Item(58, 1).Value = 'russell compression pants'
$ws.Cells.Item(59, 1).Value = 'blue knee pads for basketball'
$ws.Cells.Item(60, 1).Value = 'elbow knee pad'
$ws.Cells.Item(61, 1).Value = 'tights mens'
$ws.Cells.Item(62, 1).Value = 'super compression leggings'
$ws.Cells.Item(63, 1).Value = 'mens leggins'
$ws.Cells.Item(64, 1).Value = 'knee pads for teens'
$ws.Cells.Item(65, 1).Value = 'green mens compression pants'
$ws.Cells.Item(66, 1).Value = 'pants with padded knees'
$ws.Cells.Item(67, 1).Value = 'compression with pads'
$ws.Cells.Item(68, 1).Value = 'knee pads for sleeping'
$ws.Cells.Item(69, 1).Value = 'mens winter compression pants'
$ws.Cells.Item(70, 1).Value = 'knee pads for hvac'
$ws.Cells.Item(71, 1).Value = 'yoga pants with knee pads'
$ws.Cells.Item(72, 1).Value = 'black pants with knee pads'
$ws.Cells.Item(73, 1).Value = 'kids compression knee pads'
$ws.Cells.Item(74, 1).Value = 'elite basketball pants'
$ws.Cells.Item(75, 1).Value = 'nike leggings mens'
$ws.Cells.Item(76, 1).Value = 'compression pants men 3 pack'
$ws.Cells.Item(77, 1).Value = 'compression pants baseball'
$ws.Cells.Item(78, 1).Value = 'colored compression leggings'
$ws.Cells.Item(79, 1).Value = 'mens compression pants with pockets'
$ws.Cells.Item(80, 1).Value = 'mens compression pants xxl'
$ws.Cells.Item(81, 1).Value = 'knee pads wheels'
$ws.Cells.Item(82, 1).Value = 'ua basketball knee pads'
$ws.Cells.Item(83, 1).Value = 'protective knee pad'
$ws.Cells.Item(84, 1).Value = 'grey knee pads basketball'
$ws.Cells.Item(85, 1).Value = 'orange knee pads for basketball'
$ws.Cells.Item(86, 1).Value = 'adidas youth compression pants'
$ws.Cells.Item(87, 1).Value = 'copper compression tights for men'
$ws.Cells.Item(88, 1).Value = 'basketball knee pads youth boys mcdavid'
$ws.Cells.Item(89, 1).Value = 'nike youth basketball knee pads'
$ws.Cells.Item(90, 1).Value = 'nike kneepads'
$ws.Cells.Item(91, 1).Value = 'knee pad strap'
$ws.Cells.Item(92, 1).Value = 'internal knee pad'
$ws.Cells.Item(93, 1).Value = 'basketball legings'
$ws.Cells.Item(94, 1).Value = 'basketball knee sleves'
$ws.Cells.Item(95, 1).Value = 'under armour baseball pants men'
$ws.Cells.Item(96, 1).Value = 'cold gear compression leggings men'
$ws.Cells.Item(97, 1).Value = 'mcgregor knee pads'
$ws.Cells.Item(98, 1).Value = 'knee compression pants men'
$ws.Cells.Item(99, 1).Value = 'compression knee pads for basketball for kids'
$ws.Cells.Item(100, 1).Value = 'tesla compression pant'
